$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains two new sub-entries for "Z15" (Bodendegradation /
# Nachhaltige Landnutzung) that need to be inserted right after the
# existing Z15 rows (rows 35-36), before the Z16 rows. This shifts the
# former rows 37-42 (Z16_B01 .. Z17_B03) down by two, to rows 39-44.
#
# We shift the data "manually" (copy cell-by-cell from bottom to top)
# instead of using Rows.Insert(), and then fill in the two freshly
# freed-up rows with the new content. This keeps the existing style
# index (used by every data row) intact instead of Excel registering a
# redundant new style, which a plain row-insert operation would do.

# Move old rows 37-42 down to rows 39-44 (process bottom-up so we never
# overwrite a row before it has been read).
for ($r = 42; $r -ge 37; $r--) {
    $newRow = $r + 2
    $ws.Range("A$newRow").Value = $ws.Range("A$r").Value2
    $ws.Range("B$newRow").Value = $ws.Range("B$r").Value2
    $ws.Range("C$newRow").Value = $ws.Range("C$r").Value2
    $ws.Range("D$newRow").Value = $ws.Range("D$r").Value2
}

# Fill the newly vacated rows 37 and 38 with the two new Z15 sub-entries.
$ws.Range("A37").Value = "Z15_B03"
$ws.Range("B37").Value = "Z15"
$ws.Range("C37").Value = "Bodendegradation"
$ws.Range("D37").Value = "Bodendegradation"

$ws.Range("A38").Value = "Z15_B04"
$ws.Range("B38").Value = "Z15"
$ws.Range("C38").Value = "Nachhaltige Landnutzung"
$ws.Range("D38").Value = "XXXNachhaltige Landnutzung"

# Rows 43 and 44 are brand new rows beyond the previous data range, so
# they need the same formatting the other data rows use; copy it over
# from an existing, already-correctly-styled data row (row 36).
$ws.Range("A36:D36").Copy()
$ws.Range("A43:D44").PasteSpecial(-4122)
$excel.CutCopyMode = 0
